$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the first sheet from "skycriesmaryJSON" to "file title JSON"
# ------------------------------------------------------------------
$ws = $wb.Sheets.Item(1)
$ws.Name = "file title JSON"

# ------------------------------------------------------------------
# 2. Append the new "Star Wars stuffed animal war" rows (71-83), row 70
#    is intentionally left blank (matches a gap in the source data).
# ------------------------------------------------------------------
$items = @(
    @{ Row = 71; File = "gamemedia/cats.png";          Title = "Cat" },
    @{ Row = 72; File = "gamemedia/chickens.png";       Title = "Chicken" },
    @{ Row = 73; File = "gamemedia/crocodiles.png";     Title = "Crocodile" },
    @{ Row = 74; File = "gamemedia/darthvader.png";     Title = "Darth Vader" },
    @{ Row = 75; File = "gamemedia/dogs.png";           Title = "Dog" },
    @{ Row = 76; File = "gamemedia/lamblambs.png";      Title = "Lamb" },
    @{ Row = 77; File = "gamemedia/lions.png";          Title = "Lion" },
    @{ Row = 78; File = "gamemedia/queen.png";          Title = "Queen" },
    @{ Row = 79; File = "gamemedia/r2d2.png";           Title = "R2D2" },
    @{ Row = 80; File = "gamemedia/snakesonaplane.png"; Title = "I'm tired of these motherfuckin' snakes on this motherfuckin' plane!" },
    @{ Row = 81; File = "gamemedia/stormtrooper.png";   Title = "Storm Trooper" },
    @{ Row = 82; File = "gamemedia/thechosenone.png";   Title = "You were the chosen one!" },
    @{ Row = 83; File = "gamemedia/yoda.png";            Title = "Yoda" }
)

# Shared-string ids are minted in write order, and the source workbook
# minted all 13 titles (column B) before the 13 file names (column A) -
# replicate that exact order so the underlying sharedStrings indices line
# up, then come back through for the formulas/placeholder cells.
foreach ($item in $items) {
    $bCell = $ws.Cells.Item($item.Row, 2)
    $bCell.Value = $item.Title
    $bCell.NumberFormat = "0%"
}

foreach ($item in $items) {
    $ws.Cells.Item($item.Row, 1).Value = $item.File
}

foreach ($item in $items) {
    $r = $item.Row

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.FormulaR1C1 = "=""{'file':'""&RC[-2]&""','title':'""&RC[-1]&""'},"""

    # Rows 72-83 carry a formatted-but-empty placeholder cell in column H
    # (time format), mirroring the existing column E placeholders above.
    if ($r -ge 72) {
        $hCell = $ws.Cells.Item($r, 8)
        $hCell.Value = ""
        $hCell.NumberFormat = "h:mm"
    }
}

# ------------------------------------------------------------------
# 3. Update the view state: scroll position + active selection.
# ------------------------------------------------------------------
$aw = $excel.ActiveWindow
$aw.ScrollRow = 59
$aw.ScrollColumn = 1
$ws.Range("A83").Select()
